# "Turning off HillSlope actions as selectabie for now. Needs HDSR support."
#
# Adds a new "HillslopeDistance" data column (O) to the PlanningUnits sheet,
# and leaves the Gullies sheet/selection state the way the authored commit
# left it (PlanningUnits becomes the active/selected tab).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PlanningUnits")
$ws2 = $wb.Worksheets.Item("Gullies")

# --- New column: HillslopeDistance -----------------------------------
$ws1.Range("O1").Value = "HillslopeDistance"
$ws1.Range("O2").Value = 50817.812252999996
$ws1.Range("O3").Value = 101256.110996
$ws1.Range("O4").Value = 67849.539772999997
$ws1.Range("O5").Value = 31142.216537
$ws1.Range("O6").Value = 63657.885459999998
$ws1.Range("O7").Value = 66278.723459999994
$ws1.Range("O8").Value = 23004.958569999999

# --- Column widths (best-effort autofit to the new/longer headers) ----
$ws1.Columns.Item(9).ColumnWidth  = 29.736979166666668
$ws1.Columns.Item(10).ColumnWidth = 17.166666666666668
$ws1.Columns.Item(11).ColumnWidth = 17.451822916666668
$ws1.Columns.Item(13).ColumnWidth = 30.307291666666668
$ws1.Columns.Item(14).ColumnWidth = 12.307291666666666
$ws1.Columns.Item(15).ColumnWidth = 15.736979166666666

$ws2.Columns.Item(1).ColumnWidth = 8.592447916666666
$ws2.Columns.Item(2).ColumnWidth = 21.592447916666668
$ws2.Columns.Item(3).ColumnWidth = 11.736979166666666
$ws2.Columns.Item(4).ColumnWidth = 17.592447916666668

# --- Sheet selection / active tab -------------------------------------
# Gullies keeps a selection, but PlanningUnits becomes the active tab.
[void]$ws2.Range("F5").Select()
[void]$ws1.Activate()
[void]$ws1.Range("M13").Select()
